$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the two new test-case rows (IPA22 / IPA23) --------------------
# Order of first-use matters for shared-string layout, so write the brand
# new strings before flipping the existing "Y" results to "N".

$ws.Cells.Item(23, 1).Value = "IPA22"
$ws.Cells.Item(24, 1).Value = "IPA23"

$ws.Cells.Item(23, 2).Value = "OPQA-4853||OPQA-4854||OPQA-4856"
$ws.Cells.Item(23, 3).Value = " Verify that system provides pin option to visualization tabs when user visits dashboard Page|| Verify that system doesn't provides pin option for patents tab when user visits dashboard Page|| Verify that system dispalys all visualization when user selects pin option in every  visualizations tabs"

$ws.Cells.Item(24, 2).Value = "OPQA-4900||OPQA-4901||OPQA-4902"
$ws.Cells.Item(24, 3).Value = " Verify that system provides pin option to visualization tabs when user visits dashboard Page|| Verify that system doesn't provides pin option for patents tab when user visits dashboard Page|| Verify that system dispalys all visualization when user select"

$ws.Cells.Item(23, 4).Value = "Y"
$ws.Cells.Item(24, 4).Value = "Y"
$ws.Cells.Item(23, 5).Value = ""
$ws.Cells.Item(24, 5).Value = ""

# Match formatting of the other data rows (wrap text, border) and row heights
$ws.Range("A23:E23").Style = $ws.Range("A20:E20").Style
$ws.Range("A24:E24").Style = $ws.Range("A20:E20").Style
$ws.Rows.Item(23).RowHeight = 75
$ws.Rows.Item(24).RowHeight = 60

# --- Flip the existing results column from Y to N -----------------------
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 4).Value = "N"
}

# --- Update the view so row 19 onward + D25 are in focus ----------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("D25").Select()
